# "updated login and register"
# Append two new rows (90 and 91) to the end of the Users sheet that mirror
# the existing "register placeholder" rows (e.g. row 68): Username="moses",
# Password="bro", ID="1234", Email="m@g.c", Gender="Male", balance=0.
#
# We build the new rows by copying an existing identical row (row 68) and
# pasting it onto rows 90/91. Using Copy/PasteSpecial (instead of directly
# assigning .Value) preserves the original cell data type for the numeric-
# looking "1234" ID (kept as text, matching the source data) as well as the
# default cell style, instead of Excel's usual literal-entry auto-conversion
# turning it into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A68:F68")

$source.Copy()
$ws.Range("A90:F90").PasteSpecial(-4104)

$source.Copy()
$ws.Range("A91:F91").PasteSpecial(-4104)
